$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph and the "(c) 2020 ..." footer
# paragraph by scanning the paragraph collection for their text.
$copyrightPrefix = [string][char]169 + " 2020"

$jupiterIndex = 0
$copyrightIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "Ver no Jupiter*") {
        $jupiterIndex = $i
    }
    if ($t -like ($copyrightPrefix + "*")) {
        $copyrightIndex = $i
    }
}

if ($jupiterIndex -gt 0 -and $copyrightIndex -gt 0) {
    # Also remove the blank paragraph immediately preceding the
    # "Ver no Jupiter ..." paragraph, so the result keeps just a single
    # blank paragraph before the trailing page-break paragraph.
    $startPara = $d.Paragraphs.Item($jupiterIndex - 1)
    $endPara = $d.Paragraphs.Item($copyrightIndex)

    $delStart = $startPara.Range.Start
    $delEnd = $endPara.Range.End

    $r = $d.Range($delStart, $delEnd)
    $r.Delete()
}
